$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look like plain numbers,
# so Excel stores them verbatim as text (matching the source inline strings)
# instead of converting them to floating point numbers.
$textCells = @(
    "D5",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D13",
    "D15",
    "D16",
    "D18",
    "D19",
    "D20",
    "D22",
    "D24",
    "D25",
    "D26",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '28.194.23'
$ws.Range("E2").Value = '  -0.25%  '
$ws.Range("D3").Value = '1.911.16'
$ws.Range("E3").Value = '  +1.95%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '314.80'
$ws.Range("E5").Value = '  +0.84%  '
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("D7").Value = '0.5070'
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").Value = '0.3919'
$ws.Range("E8").Value = '  -0.54%  '
$ws.Range("D9").Value = '0.09277'
$ws.Range("E9").Value = '  -3.61%  '
$ws.Range("D10").Value = '1.140'
$ws.Range("E10").Value = '  -0.51%  '
$ws.Range("D11").Value = '41.85'
$ws.Range("E11").Value = '  +2.39%  '
$ws.Range("D12").Value = '6.387'
$ws.Range("E12").Value = '  -1.52%  '
$ws.Range("D13").Value = '20.91'
$ws.Range("E13").Value = '  -0.49%  '
$ws.Range("D14").Value = '1.907.30'
$ws.Range("E14").Value = '  +1.32%  '
$ws.Range("D15").Value = '7.310'
$ws.Range("E15").Value = '  -1.67%  '
$ws.Range("D16").Value = '1.002'
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("E17").Value = '  -1.01%  '
$ws.Range("D18").Value = '92.38'
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("D19").Value = '0.06604'
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").Value = '17.98'
$ws.Range("E20").Value = '  +2.12%  '
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("D22").Value = '6.224'
$ws.Range("E22").Value = '  +0.67%  '
$ws.Range("D23").Value = '28.256.07'
$ws.Range("E23").Value = '  -0.26%  '
$ws.Range("D24").Value = '11.44'
$ws.Range("E24").Value = '  +1.19%  '
$ws.Range("D25").Value = '2.322'
$ws.Range("E25").Value = '  +0.71%  '
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = '2.589'
$ws.Range("E26").Value = '  +0.72%  '
$ws.Range("B27").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C27").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D27").Value = '2.128.04'
$ws.Range("E27").Value = '  +1.46%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '21.07'
$ws.Range("E28").Value = '  -0.77%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").Value = '157.94'
$ws.Range("E29").Value = '  -0.60%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = '127.10'
$ws.Range("E30").Value = '  -0.40%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '1.098'
$ws.Range("E31").Value = '  +2.76%  '
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = '0.1075'
$ws.Range("E32").Value = '  +0.44%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '5.637'
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '3.618'
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("B35").Value = 'FraxShare'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D35").Value = '9.685'
$ws.Range("E35").Value = '  +1.37%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '0.06672'
$ws.Range("E36").Value = '  -0.56%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.02422'
$ws.Range("E37").Value = '  +1.23%  '
$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").Value = '1.247'
$ws.Range("E38").Value = '  -0.66%  '
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").Value = '0.2197'
$ws.Range("E39").Value = '  +0.13%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '1.290'
$ws.Range("E40").Value = '  +8.64%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.6441'
$ws.Range("E41").Value = '  +0.94%  '
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").Value = '11.49'
$ws.Range("E42").Value = '  -0.24%  '
$ws.Range("B43").Value = 'InternetComputer(DFINITY)'
$ws.Range("C43").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D43").Value = '4.990'
$ws.Range("E43").Value = '  -0.23%  '
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  -0.16%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '13.27'
$ws.Range("E45").Value = '  -1.46%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.6031'
$ws.Range("E46").Value = '  +0.62%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").Value = '3.723'
$ws.Range("E47").Value = '  +1.68%  '
$ws.Range("B48").Value = 'WEMIXTOKEN'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = '1.280'
$ws.Range("E48").Value = '  +0.92%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '2.021'
$ws.Range("E49").Value = '  +0.70%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '122.88'
$ws.Range("E50").Value = '  -1.08%  '
$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").Value = '1.185'
$ws.Range("E51").Value = '  -0.87%  '
